# Insert two new weekly price records for "Macroferia Regional de Talca - Choclo"
# right above the existing row 239, shifting the existing 239:322 block down to
# 241:324 (matches dimension growing from A1:R322 to A1:R324).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("239:240").Insert()

# New row 239: Choclero / Primera
$ws.Cells.Item(239, 1).Value = 5
$ws.Cells.Item(239, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(239, 3).Value = "Maule"
$ws.Cells.Item(239, 4).Value = [DateTime]"2023-03-03"
$ws.Cells.Item(239, 5).Value = 7
$ws.Cells.Item(239, 6).Value = 100112024
$ws.Cells.Item(239, 7).Value = "Choclo"
$ws.Cells.Item(239, 8).Value = "Choclero"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 20000
$ws.Cells.Item(239, 11).Value = 500
$ws.Cells.Item(239, 12).Value = 500
$ws.Cells.Item(239, 13).Value = 500
$ws.Cells.Item(239, 14).Value = "`$/unidad"
$ws.Cells.Item(239, 15).Value = "Región del Maule"
$ws.Cells.Item(239, 16).Value = 500
$ws.Cells.Item(239, 17).Value = 1
$ws.Cells.Item(239, 18).Value = "Hortaliza"

# New row 240: Choclero / Segunda
$ws.Cells.Item(240, 1).Value = 5
$ws.Cells.Item(240, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(240, 3).Value = "Maule"
$ws.Cells.Item(240, 4).Value = [DateTime]"2023-03-03"
$ws.Cells.Item(240, 5).Value = 7
$ws.Cells.Item(240, 6).Value = 100112024
$ws.Cells.Item(240, 7).Value = "Choclo"
$ws.Cells.Item(240, 8).Value = "Choclero"
$ws.Cells.Item(240, 9).Value = "Segunda"
$ws.Cells.Item(240, 10).Value = 20000
$ws.Cells.Item(240, 11).Value = 400
$ws.Cells.Item(240, 12).Value = 400
$ws.Cells.Item(240, 13).Value = 400
$ws.Cells.Item(240, 14).Value = "`$/unidad"
$ws.Cells.Item(240, 15).Value = "Región del Maule"
$ws.Cells.Item(240, 16).Value = 400
$ws.Cells.Item(240, 17).Value = 1
$ws.Cells.Item(240, 18).Value = "Hortaliza"
